$d = $word.ActiveDocument

# Update the title date line
$d.Content.Find.Execute("2026-01-09 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2026-01-10 Saturday", 2)

# Update the division-problem answer table.
# The table has 20 rows x 5 columns; data lives in rows 1, 5, 9, 13, 17
# (the other rows are blank spacer rows). Cells are updated by explicit
# (row, column) address so duplicate source strings do not collide.
$t = $d.Tables.Item(1)

$updates = @(
    @{ r = 1;  c = 1; v = "57÷7=8, 1" },
    @{ r = 1;  c = 2; v = "20÷7=2, 6" },
    @{ r = 1;  c = 3; v = "34÷5=6, 4" },
    @{ r = 1;  c = 4; v = "91÷9=10, 1" },
    @{ r = 1;  c = 5; v = "14÷8=1, 6" },

    @{ r = 5;  c = 1; v = "36÷2=18, 0" },
    @{ r = 5;  c = 2; v = "49÷5=9, 4" },
    @{ r = 5;  c = 3; v = "96÷4=24, 0" },
    @{ r = 5;  c = 4; v = "53÷8=6, 5" },
    @{ r = 5;  c = 5; v = "67÷6=11, 1" },

    @{ r = 9;  c = 1; v = "16÷4=4, 0" },
    @{ r = 9;  c = 2; v = "63÷6=10, 3" },
    @{ r = 9;  c = 3; v = "35÷4=8, 3" },
    @{ r = 9;  c = 4; v = "83÷8=10, 3" },
    @{ r = 9;  c = 5; v = "82÷4=20, 2" },

    @{ r = 13; c = 1; v = "20÷4=5, 0" },
    @{ r = 13; c = 2; v = "73÷8=9, 1" },
    @{ r = 13; c = 3; v = "94÷5=18, 4" },
    @{ r = 13; c = 4; v = "20÷6=3, 2" },
    @{ r = 13; c = 5; v = "86÷4=21, 2" },

    @{ r = 17; c = 1; v = "73÷3=24, 1" },
    @{ r = 17; c = 2; v = "44÷7=6, 2" },
    @{ r = 17; c = 3; v = "72÷9=8, 0" },
    @{ r = 17; c = 4; v = "50÷4=12, 2" },
    @{ r = 17; c = 5; v = "97÷7=13, 6" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.r, $u.c)
    $cell.Range.Text = $u.v
}
